$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the "Overall Demand" row (row 2) values for columns B, C, D.
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 544.9647926184913
$ws.Range("D2").Value = 1057.634047914892
